$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 changes
$ws.Range("M2").Value = 1.11
$ws.Range("N2").Value = 6.5

# Row 3 changes
$ws.Range("BD3").Value = 151

# Row 4 changes
$ws.Range("G4").Value = 1.95
$ws.Range("I4").Value = 4.75
$ws.Range("L4").Value = 6
$ws.Range("O4").Value = 1.67
$ws.Range("P4").Value = 2.1
$ws.Range("X4").Value = 7
$ws.Range("AE4").Value = 29
$ws.Range("AG4").Value = 8
$ws.Range("AO4").Value = 12
$ws.Range("AQ4").Value = 41
$ws.Range("AU4").Value = 12
$ws.Range("AZ4").Value = 151
